$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.2335874793234083
$ws.Range("E2").Value = 0.2335874793234083

# Row 3
$ws.Range("D3").Value = 0.3173327494528142
$ws.Range("E3").Value = 0.3173327494528142

# Row 4
$ws.Range("D4").Value = 0.1314683145534435
$ws.Range("E4").Value = 0.1314683145534435

# Row 5
$ws.Range("D5").Value = 0.3076367367788183
$ws.Range("E5").Value = 0.3076367367788183

# Row 6
$ws.Range("D6").Value = 0.2809949785630151
$ws.Range("E6").Value = 0.2809949785630151

# Row 7
$ws.Range("D7").Value = 0.4550958698018167
$ws.Range("E7").Value = 0.5449041301981834

# Row 8
$ws.Range("D8").Value = 0.02178231699227287
$ws.Range("E8").Value = 0.9782176830077272

# Row 9
$ws.Range("D9").Value = 0.4587394875733926
$ws.Range("E9").Value = 0.5412605124266074

# Row 10
$ws.Range("C10").Value = $true
$ws.Range("D10").Value = 0.6481054049809666
$ws.Range("E10").Value = 0.3518945950190334

# Row 11
$ws.Range("C11").Value = $true
$ws.Range("D11").Value = 0.5450179765751472
$ws.Range("E11").Value = 0.4549820234248528
$ws.Range("F11").Value = 0.7920080423355103
$ws.Range("G11").Value = 0.7

# Row 12
$ws.Range("D12").Value = 0.2084977903266872
$ws.Range("E12").Value = 0.2084977903266872

# Row 13
$ws.Range("D13").Value = 0.1957666513813293
$ws.Range("E13").Value = 0.1957666513813293

# Row 14
$ws.Range("D14").Value = 0.04151338055776142
$ws.Range("E14").Value = 0.04151338055776142

# Row 15
$ws.Range("D15").Value = 0.1936490800145129
$ws.Range("E15").Value = 0.1936490800145129

# Row 16
$ws.Range("D16").Value = 0.1950314156639996
$ws.Range("E16").Value = 0.1950314156639996

# Row 17
$ws.Range("D17").Value = 0.3711061322852943
$ws.Range("E17").Value = 0.6288938677147058

# Row 18
$ws.Range("D18").Value = 0.002842995850197117
$ws.Range("E18").Value = 0.9971570041498029

# Row 19
$ws.Range("D19").Value = 0.3740249317764517
$ws.Range("E19").Value = 0.6259750682235483

# Row 20
$ws.Range("C20").Value = $true
$ws.Range("D20").Value = 0.6178413825779572
$ws.Range("E20").Value = 0.3821586174220428

# Row 21
$ws.Range("C21").Value = $true
$ws.Range("D21").Value = 0.4464085480326587
$ws.Range("E21").Value = 0.5535914519673413
$ws.Range("F21").Value = 1.005191802978516
$ws.Range("G21").Value = 0.6

$wb.Save()
